# pontos notáveis - incremento na tabela de ranking
# Multiply the values in columns E and F (rows 2-7) by 100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..7) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value2 = $cell.Value2 * 100
    }
}
